$wb = $excel.ActiveWorkbook

# --- Proximity sheet: add rows 14-21 ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-02-01"
$ws.Cells.Item(14,2).Value = "16:07:39"
$ws.Cells.Item(14,3).Value = "16:00"
$ws.Cells.Item(14,4).Value = "Living Room Main Door"
$ws.Cells.Item(14,5).Value = "ENTER"
$ws.Cells.Item(14,6).Value = "User ENTERED Living Room Main Door"

$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-02-01"
$ws.Cells.Item(15,2).Value = "16:07:45"
$ws.Cells.Item(15,3).Value = "16:00"
$ws.Cells.Item(15,4).Value = "Living Room Main Door"
$ws.Cells.Item(15,5).Value = "EXIT"
$ws.Cells.Item(15,6).Value = "User EXITED Living Room Main Door"

$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-02-01"
$ws.Cells.Item(16,2).Value = "16:07:52"
$ws.Cells.Item(16,3).Value = "16:00"
$ws.Cells.Item(16,4).Value = "Living Room Main Door"
$ws.Cells.Item(16,5).Value = "ENTER"
$ws.Cells.Item(16,6).Value = "User ENTERED Living Room Main Door"

$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2026-02-01"
$ws.Cells.Item(17,2).Value = "16:07:58"
$ws.Cells.Item(17,3).Value = "16:00"
$ws.Cells.Item(17,4).Value = "Living Room Main Door"
$ws.Cells.Item(17,5).Value = "EXIT"
$ws.Cells.Item(17,6).Value = "User EXITED Living Room Main Door"

$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2026-02-01"
$ws.Cells.Item(18,2).Value = "16:08:17"
$ws.Cells.Item(18,3).Value = "16:00"
$ws.Cells.Item(18,4).Value = "Living Room Main Door"
$ws.Cells.Item(18,5).Value = "ENTER"
$ws.Cells.Item(18,6).Value = "User ENTERED Living Room Main Door"

$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = "2026-02-01"
$ws.Cells.Item(19,2).Value = "16:08:19"
$ws.Cells.Item(19,3).Value = "16:00"
$ws.Cells.Item(19,4).Value = "Living Room Main Door"
$ws.Cells.Item(19,5).Value = "EXIT"
$ws.Cells.Item(19,6).Value = "User EXITED Living Room Main Door"

$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = "2026-02-01"
$ws.Cells.Item(20,2).Value = "16:08:23"
$ws.Cells.Item(20,3).Value = "16:00"
$ws.Cells.Item(20,4).Value = "Living Room Main Door"
$ws.Cells.Item(20,5).Value = "EXIT"
$ws.Cells.Item(20,6).Value = "User EXITED Living Room Main Door"

$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = "2026-02-01"
$ws.Cells.Item(21,2).Value = "16:08:28"
$ws.Cells.Item(21,3).Value = "16:00"
$ws.Cells.Item(21,4).Value = "Living Room Main Door"
$ws.Cells.Item(21,5).Value = "ENTER"
$ws.Cells.Item(21,6).Value = "User ENTERED Living Room Main Door"

# --- Camera sheet: add rows 13-16 ---
$ws = $wb.Worksheets.Item("Camera")
$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = "2026-02-01"
$ws.Cells.Item(13,2).Value = "16:07:41"
$ws.Cells.Item(13,3).Value = "16:00"
$ws.Cells.Item(13,4).Value = "Living Room Main Door"
$ws.Cells.Item(13,5).Value = "Image Captured"
$ws.Cells.Item(13,6).Value = "Active"

$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-02-01"
$ws.Cells.Item(14,2).Value = "16:07:53"
$ws.Cells.Item(14,3).Value = "16:00"
$ws.Cells.Item(14,4).Value = "Living Room Main Door"
$ws.Cells.Item(14,5).Value = "Image Captured"
$ws.Cells.Item(14,6).Value = "Active"

$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-02-01"
$ws.Cells.Item(15,2).Value = "16:08:19"
$ws.Cells.Item(15,3).Value = "16:00"
$ws.Cells.Item(15,4).Value = "Living Room Main Door"
$ws.Cells.Item(15,5).Value = "Image Captured"
$ws.Cells.Item(15,6).Value = "Active"

$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-02-01"
$ws.Cells.Item(16,2).Value = "16:08:30"
$ws.Cells.Item(16,3).Value = "16:00"
$ws.Cells.Item(16,4).Value = "Living Room Main Door"
$ws.Cells.Item(16,5).Value = "Image Captured"
$ws.Cells.Item(16,6).Value = "Active"

# --- mmWave(BR) sheet: add rows 2-18 (numeric breathing-rate values) ---
$ws = $wb.Worksheets.Item("mmWave(BR)")
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "2026-02-01"
$ws.Cells.Item(2,2).Value = "16:05:41"
$ws.Cells.Item(2,3).Value = "16:00"
$ws.Cells.Item(2,4).Value = "Bedroom"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = "Empty"

$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2026-02-01"
$ws.Cells.Item(3,2).Value = "16:05:42"
$ws.Cells.Item(3,3).Value = "16:00"
$ws.Cells.Item(3,4).Value = "Bedroom"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = "Occupied"

$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "2026-02-01"
$ws.Cells.Item(4,2).Value = "16:05:42"
$ws.Cells.Item(4,3).Value = "16:00"
$ws.Cells.Item(4,4).Value = "Bedroom"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = "Occupied"

$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "2026-02-01"
$ws.Cells.Item(5,2).Value = "16:05:47"
$ws.Cells.Item(5,3).Value = "16:00"
$ws.Cells.Item(5,4).Value = "Bedroom"
$ws.Cells.Item(5,5).Value = 48
$ws.Cells.Item(5,6).Value = "Occupied"

$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "2026-02-01"
$ws.Cells.Item(6,2).Value = "16:05:48"
$ws.Cells.Item(6,3).Value = "16:00"
$ws.Cells.Item(6,4).Value = "Bedroom"
$ws.Cells.Item(6,5).Value = 28
$ws.Cells.Item(6,6).Value = "Occupied"

$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "2026-02-01"
$ws.Cells.Item(7,2).Value = "16:05:49"
$ws.Cells.Item(7,3).Value = "16:00"
$ws.Cells.Item(7,4).Value = "Bedroom"
$ws.Cells.Item(7,5).Value = 34
$ws.Cells.Item(7,6).Value = "Occupied"

$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "2026-02-01"
$ws.Cells.Item(8,2).Value = "16:05:50"
$ws.Cells.Item(8,3).Value = "16:00"
$ws.Cells.Item(8,4).Value = "Bedroom"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = "Occupied"

$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = "2026-02-01"
$ws.Cells.Item(9,2).Value = "16:05:58"
$ws.Cells.Item(9,3).Value = "16:00"
$ws.Cells.Item(9,4).Value = "Bedroom"
$ws.Cells.Item(9,5).Value = 19
$ws.Cells.Item(9,6).Value = "Occupied"

$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "2026-02-01"
$ws.Cells.Item(10,2).Value = "16:05:59"
$ws.Cells.Item(10,3).Value = "16:00"
$ws.Cells.Item(10,4).Value = "Bedroom"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = "Occupied"

$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "2026-02-01"
$ws.Cells.Item(11,2).Value = "16:06:03"
$ws.Cells.Item(11,3).Value = "16:00"
$ws.Cells.Item(11,4).Value = "Bedroom"
$ws.Cells.Item(11,5).Value = 4
$ws.Cells.Item(11,6).Value = "Occupied"

$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = "2026-02-01"
$ws.Cells.Item(12,2).Value = "16:06:04"
$ws.Cells.Item(12,3).Value = "16:00"
$ws.Cells.Item(12,4).Value = "Bedroom"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = "Occupied"

$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = "2026-02-01"
$ws.Cells.Item(13,2).Value = "16:06:12"
$ws.Cells.Item(13,3).Value = "16:00"
$ws.Cells.Item(13,4).Value = "Bedroom"
$ws.Cells.Item(13,5).Value = 23
$ws.Cells.Item(13,6).Value = "Occupied"

$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-02-01"
$ws.Cells.Item(14,2).Value = "16:06:13"
$ws.Cells.Item(14,3).Value = "16:00"
$ws.Cells.Item(14,4).Value = "Bedroom"
$ws.Cells.Item(14,5).Value = 24
$ws.Cells.Item(14,6).Value = "Occupied"

$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-02-01"
$ws.Cells.Item(15,2).Value = "16:06:14"
$ws.Cells.Item(15,3).Value = "16:00"
$ws.Cells.Item(15,4).Value = "Bedroom"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = "Occupied"

$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-02-01"
$ws.Cells.Item(16,2).Value = "16:06:15"
$ws.Cells.Item(16,3).Value = "16:00"
$ws.Cells.Item(16,4).Value = "Bedroom"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = "Occupied"

$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2026-02-01"
$ws.Cells.Item(17,2).Value = "16:08:01"
$ws.Cells.Item(17,3).Value = "16:00"
$ws.Cells.Item(17,4).Value = "Bedroom"
$ws.Cells.Item(17,5).Value = 32
$ws.Cells.Item(17,6).Value = "Occupied"

$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2026-02-01"
$ws.Cells.Item(18,2).Value = "16:08:02"
$ws.Cells.Item(18,3).Value = "16:00"
$ws.Cells.Item(18,4).Value = "Bedroom"
$ws.Cells.Item(18,5).Value = 2
$ws.Cells.Item(18,6).Value = "Occupied"

# --- mmWave(HR) sheet: add rows 2-18 (numeric heart-rate values) ---
$ws = $wb.Worksheets.Item("mmWave(HR)")
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "2026-02-01"
$ws.Cells.Item(2,2).Value = "16:05:40"
$ws.Cells.Item(2,3).Value = "16:00"
$ws.Cells.Item(2,4).Value = "Bedroom"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = "Empty"

$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2026-02-01"
$ws.Cells.Item(3,2).Value = "16:05:42"
$ws.Cells.Item(3,3).Value = "16:00"
$ws.Cells.Item(3,4).Value = "Bedroom"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = "Occupied"

$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "2026-02-01"
$ws.Cells.Item(4,2).Value = "16:05:42"
$ws.Cells.Item(4,3).Value = "16:00"
$ws.Cells.Item(4,4).Value = "Bedroom"
$ws.Cells.Item(4,5).Value = 50
$ws.Cells.Item(4,6).Value = "Occupied"

$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "2026-02-01"
$ws.Cells.Item(5,2).Value = "16:05:47"
$ws.Cells.Item(5,3).Value = "16:00"
$ws.Cells.Item(5,4).Value = "Bedroom"
$ws.Cells.Item(5,5).Value = 96
$ws.Cells.Item(5,6).Value = "Occupied"

$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "2026-02-01"
$ws.Cells.Item(6,2).Value = "16:05:48"
$ws.Cells.Item(6,3).Value = "16:00"
$ws.Cells.Item(6,4).Value = "Bedroom"
$ws.Cells.Item(6,5).Value = 76
$ws.Cells.Item(6,6).Value = "Occupied"

$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "2026-02-01"
$ws.Cells.Item(7,2).Value = "16:05:49"
$ws.Cells.Item(7,3).Value = "16:00"
$ws.Cells.Item(7,4).Value = "Bedroom"
$ws.Cells.Item(7,5).Value = 82
$ws.Cells.Item(7,6).Value = "Occupied"

$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "2026-02-01"
$ws.Cells.Item(8,2).Value = "16:05:50"
$ws.Cells.Item(8,3).Value = "16:00"
$ws.Cells.Item(8,4).Value = "Bedroom"
$ws.Cells.Item(8,5).Value = 50
$ws.Cells.Item(8,6).Value = "Occupied"

$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = "2026-02-01"
$ws.Cells.Item(9,2).Value = "16:05:57"
$ws.Cells.Item(9,3).Value = "16:00"
$ws.Cells.Item(9,4).Value = "Bedroom"
$ws.Cells.Item(9,5).Value = 67
$ws.Cells.Item(9,6).Value = "Occupied"

$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "2026-02-01"
$ws.Cells.Item(10,2).Value = "16:05:58"
$ws.Cells.Item(10,3).Value = "16:00"
$ws.Cells.Item(10,4).Value = "Bedroom"
$ws.Cells.Item(10,5).Value = 50
$ws.Cells.Item(10,6).Value = "Occupied"

$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "2026-02-01"
$ws.Cells.Item(11,2).Value = "16:06:03"
$ws.Cells.Item(11,3).Value = "16:00"
$ws.Cells.Item(11,4).Value = "Bedroom"
$ws.Cells.Item(11,5).Value = 52
$ws.Cells.Item(11,6).Value = "Occupied"

$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = "2026-02-01"
$ws.Cells.Item(12,2).Value = "16:06:03"
$ws.Cells.Item(12,3).Value = "16:00"
$ws.Cells.Item(12,4).Value = "Bedroom"
$ws.Cells.Item(12,5).Value = 50
$ws.Cells.Item(12,6).Value = "Occupied"

$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = "2026-02-01"
$ws.Cells.Item(13,2).Value = "16:06:12"
$ws.Cells.Item(13,3).Value = "16:00"
$ws.Cells.Item(13,4).Value = "Bedroom"
$ws.Cells.Item(13,5).Value = 71
$ws.Cells.Item(13,6).Value = "Occupied"

$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-02-01"
$ws.Cells.Item(14,2).Value = "16:06:12"
$ws.Cells.Item(14,3).Value = "16:00"
$ws.Cells.Item(14,4).Value = "Bedroom"
$ws.Cells.Item(14,5).Value = 72
$ws.Cells.Item(14,6).Value = "Occupied"

$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-02-01"
$ws.Cells.Item(15,2).Value = "16:06:14"
$ws.Cells.Item(15,3).Value = "16:00"
$ws.Cells.Item(15,4).Value = "Bedroom"
$ws.Cells.Item(15,5).Value = 51
$ws.Cells.Item(15,6).Value = "Occupied"

$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-02-01"
$ws.Cells.Item(16,2).Value = "16:06:15"
$ws.Cells.Item(16,3).Value = "16:00"
$ws.Cells.Item(16,4).Value = "Bedroom"
$ws.Cells.Item(16,5).Value = 50
$ws.Cells.Item(16,6).Value = "Occupied"

$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2026-02-01"
$ws.Cells.Item(17,2).Value = "16:08:01"
$ws.Cells.Item(17,3).Value = "16:00"
$ws.Cells.Item(17,4).Value = "Bedroom"
$ws.Cells.Item(17,5).Value = 80
$ws.Cells.Item(17,6).Value = "Occupied"

$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2026-02-01"
$ws.Cells.Item(18,2).Value = "16:08:02"
$ws.Cells.Item(18,3).Value = "16:00"
$ws.Cells.Item(18,4).Value = "Bedroom"
$ws.Cells.Item(18,5).Value = 50
$ws.Cells.Item(18,6).Value = "Occupied"

# --- mmWave(InBed) sheet: add rows 2-18 (text in/out-of-bed state) ---
$ws = $wb.Worksheets.Item("mmWave(InBed)")
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "2026-02-01"
$ws.Cells.Item(2,2).Value = "16:05:40"
$ws.Cells.Item(2,3).Value = "16:00"
$ws.Cells.Item(2,4).Value = "Bedroom"
$ws.Cells.Item(2,5).Value = "Out of Bed"
$ws.Cells.Item(2,6).Value = "Empty"

$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2026-02-01"
$ws.Cells.Item(3,2).Value = "16:05:41"
$ws.Cells.Item(3,3).Value = "16:00"
$ws.Cells.Item(3,4).Value = "Bedroom"
$ws.Cells.Item(3,5).Value = "In Bed"
$ws.Cells.Item(3,6).Value = "Occupied"

$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "2026-02-01"
$ws.Cells.Item(4,2).Value = "16:05:42"
$ws.Cells.Item(4,3).Value = "16:00"
$ws.Cells.Item(4,4).Value = "Bedroom"
$ws.Cells.Item(4,5).Value = "In Bed"
$ws.Cells.Item(4,6).Value = "Occupied"

$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "2026-02-01"
$ws.Cells.Item(5,2).Value = "16:05:46"
$ws.Cells.Item(5,3).Value = "16:00"
$ws.Cells.Item(5,4).Value = "Bedroom"
$ws.Cells.Item(5,5).Value = "In Bed"
$ws.Cells.Item(5,6).Value = "Occupied"

$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "2026-02-01"
$ws.Cells.Item(6,2).Value = "16:05:47"
$ws.Cells.Item(6,3).Value = "16:00"
$ws.Cells.Item(6,4).Value = "Bedroom"
$ws.Cells.Item(6,5).Value = "In Bed"
$ws.Cells.Item(6,6).Value = "Occupied"

$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "2026-02-01"
$ws.Cells.Item(7,2).Value = "16:05:48"
$ws.Cells.Item(7,3).Value = "16:00"
$ws.Cells.Item(7,4).Value = "Bedroom"
$ws.Cells.Item(7,5).Value = "In Bed"
$ws.Cells.Item(7,6).Value = "Occupied"

$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "2026-02-01"
$ws.Cells.Item(8,2).Value = "16:05:49"
$ws.Cells.Item(8,3).Value = "16:00"
$ws.Cells.Item(8,4).Value = "Bedroom"
$ws.Cells.Item(8,5).Value = "In Bed"
$ws.Cells.Item(8,6).Value = "Occupied"

$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = "2026-02-01"
$ws.Cells.Item(9,2).Value = "16:05:57"
$ws.Cells.Item(9,3).Value = "16:00"
$ws.Cells.Item(9,4).Value = "Bedroom"
$ws.Cells.Item(9,5).Value = "In Bed"
$ws.Cells.Item(9,6).Value = "Occupied"

$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = "2026-02-01"
$ws.Cells.Item(10,2).Value = "16:05:58"
$ws.Cells.Item(10,3).Value = "16:00"
$ws.Cells.Item(10,4).Value = "Bedroom"
$ws.Cells.Item(10,5).Value = "In Bed"
$ws.Cells.Item(10,6).Value = "Occupied"

$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "2026-02-01"
$ws.Cells.Item(11,2).Value = "16:06:03"
$ws.Cells.Item(11,3).Value = "16:00"
$ws.Cells.Item(11,4).Value = "Bedroom"
$ws.Cells.Item(11,5).Value = "In Bed"
$ws.Cells.Item(11,6).Value = "Occupied"

$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = "2026-02-01"
$ws.Cells.Item(12,2).Value = "16:06:03"
$ws.Cells.Item(12,3).Value = "16:00"
$ws.Cells.Item(12,4).Value = "Bedroom"
$ws.Cells.Item(12,5).Value = "In Bed"
$ws.Cells.Item(12,6).Value = "Occupied"

$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = "2026-02-01"
$ws.Cells.Item(13,2).Value = "16:06:12"
$ws.Cells.Item(13,3).Value = "16:00"
$ws.Cells.Item(13,4).Value = "Bedroom"
$ws.Cells.Item(13,5).Value = "In Bed"
$ws.Cells.Item(13,6).Value = "Occupied"

$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-02-01"
$ws.Cells.Item(14,2).Value = "16:06:12"
$ws.Cells.Item(14,3).Value = "16:00"
$ws.Cells.Item(14,4).Value = "Bedroom"
$ws.Cells.Item(14,5).Value = "In Bed"
$ws.Cells.Item(14,6).Value = "Occupied"

$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-02-01"
$ws.Cells.Item(15,2).Value = "16:06:13"
$ws.Cells.Item(15,3).Value = "16:00"
$ws.Cells.Item(15,4).Value = "Bedroom"
$ws.Cells.Item(15,5).Value = "In Bed"
$ws.Cells.Item(15,6).Value = "Occupied"

$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-02-01"
$ws.Cells.Item(16,2).Value = "16:06:15"
$ws.Cells.Item(16,3).Value = "16:00"
$ws.Cells.Item(16,4).Value = "Bedroom"
$ws.Cells.Item(16,5).Value = "In Bed"
$ws.Cells.Item(16,6).Value = "Occupied"

$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2026-02-01"
$ws.Cells.Item(17,2).Value = "16:08:00"
$ws.Cells.Item(17,3).Value = "16:00"
$ws.Cells.Item(17,4).Value = "Bedroom"
$ws.Cells.Item(17,5).Value = "In Bed"
$ws.Cells.Item(17,6).Value = "Occupied"

$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2026-02-01"
$ws.Cells.Item(18,2).Value = "16:08:01"
$ws.Cells.Item(18,3).Value = "16:00"
$ws.Cells.Item(18,4).Value = "Bedroom"
$ws.Cells.Item(18,5).Value = "In Bed"
$ws.Cells.Item(18,6).Value = "Occupied"

